$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Table2 currently spans A1:E3 (header + 2 data rows). Grow it by one row
# via the ListObject so the table ref/autoFilter and sheet dimension follow.
$tbl = $ws.ListObjects.Item(1)
$newListRow = $tbl.ListRows.Add()

# New row (row 4) values.
$ws.Range("A4").Value = "1527. Patients With a Condition"
$ws.Range("B4").Value = "Easy"
$ws.Range("C4").Value = "String Methods"
$ws.Range("D4").Value = "Use str.contains() with RegEx"

# B4 ("Easy") reuses the same green-fill format as B2/B3.
$ws.Range("B2").Copy()
$ws.Range("B4").PasteSpecial(-4122)
$ws.Range("B4").Value = "Easy"
$excel.CutCopyMode = $false

$linkCell = $ws.Range("E4")
$ws.Hyperlinks.Add($linkCell, "https://leetcode.com/problems/patients-with-a-condition/solutions/3853981/pandas-mysql-an-effortless-and-simple-approach-with-comments-and-explanation/?envType=study-plan-v2&envId=30-days-of-pandas&lang=pythondata ")
$linkCell.Style = "Hyperlink"

# Match the saved selection state in the diff.
$ws.Range("E13").Select()
